# Fruta / hortaliza, semanal
# Rotate the weekly price records: row2 <- old row3, row3 <- old row4, row4 <- old row2
# (columns D, L, M, N, O, P, R, S carry the weekly-varying data)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original values before overwriting anything
$D2 = $ws.Range("D2").Value2
$L2 = $ws.Range("L2").Value2
$M2 = $ws.Range("M2").Value2
$N2 = $ws.Range("N2").Value2
$O2 = $ws.Range("O2").Value2
$P2 = $ws.Range("P2").Value2
$R2 = $ws.Range("R2").Value2
$S2 = $ws.Range("S2").Value2

$D3 = $ws.Range("D3").Value2
$L3 = $ws.Range("L3").Value2
$M3 = $ws.Range("M3").Value2
$N3 = $ws.Range("N3").Value2
$O3 = $ws.Range("O3").Value2
$P3 = $ws.Range("P3").Value2
$R3 = $ws.Range("R3").Value2
$S3 = $ws.Range("S3").Value2

$D4 = $ws.Range("D4").Value2
$L4 = $ws.Range("L4").Value2
$M4 = $ws.Range("M4").Value2
$N4 = $ws.Range("N4").Value2
$O4 = $ws.Range("O4").Value2
$P4 = $ws.Range("P4").Value2
$R4 = $ws.Range("R4").Value2
$S4 = $ws.Range("S4").Value2

# Row 2 <- old Row 3
$ws.Range("D2").Value2 = $D3
$ws.Range("L2").Value2 = $L3
$ws.Range("M2").Value2 = $M3
$ws.Range("N2").Value2 = $N3
$ws.Range("O2").Value2 = $O3
$ws.Range("P2").Value2 = $P3
$ws.Range("R2").Value2 = $R3
$ws.Range("S2").Value2 = $S3

# Row 3 <- old Row 4
$ws.Range("D3").Value2 = $D4
$ws.Range("L3").Value2 = $L4
$ws.Range("M3").Value2 = $M4
$ws.Range("N3").Value2 = $N4
$ws.Range("O3").Value2 = $O4
$ws.Range("P3").Value2 = $P4
$ws.Range("R3").Value2 = $R4
$ws.Range("S3").Value2 = $S4

# Row 4 <- old Row 2
$ws.Range("D4").Value2 = $D2
$ws.Range("L4").Value2 = $L2
$ws.Range("M4").Value2 = $M2
$ws.Range("N4").Value2 = $N2
$ws.Range("O4").Value2 = $O2
$ws.Range("P4").Value2 = $P2
$ws.Range("R4").Value2 = $R2
$ws.Range("S4").Value2 = $S2
